$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns (Month, Day, Year) right before the existing
# "Date Sampled" column (old E), shifting everything from E onward to H onward.
$ws.Range("E1:G1").EntireColumn.Insert() | Out-Null

# Match the width of the neighboring "bestFit" text columns (C:D, width 9.5)
$ws.Range("E1:G1").ColumnWidth = 8.67

# Header labels for the new columns
$ws.Range("E1").Value = "Month"
$ws.Range("F1").Value = "Day"
$ws.Range("G1").Value = "Year"

# Row 2 (Transect 1): Date Sampled 8/2/2016
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = 2016

# Row 3 (Transect 2): Date Sampled 8/2/2016
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 2016

# Row 4 (Transect 3): Date Sampled 8/4/2016
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 4
$ws.Range("G4").Value = 2016

# Row 5 (Transect 4): Date Sampled 8/4/2016
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 4
$ws.Range("G5").Value = 2016

# Match the final selected cell left behind in the authored file
$ws.Range("G5").Select() | Out-Null

$wb.Save()
